# This script reassigns several bird-observation records (rows 11-23, sheet
# "Artfynd") to different rows, updating the identifying fields (Id,
# taxon/species info, activity, coordinates, accuracy, comment, observer)
# for each row so that the sighting records end up on the rows described by
# the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 ---
$ws.Range("A11").Value = 107357712
$ws.Range("I11").Value = $null
$ws.Range("M11").Value = "äldre spår"
$ws.Range("Q11").Value = 421926.4665079587
$ws.Range("R11").Value = 7050030.281593022
$ws.Range("AC11").Value = "ringhack"

# --- Row 12 ---
$ws.Range("A12").Value = 107357711
$ws.Range("M12").Value = "färska spår"
$ws.Range("Q12").Value = 421929.8157829407
$ws.Range("R12").Value = 7050021.285223413

# --- Row 13 ---
$ws.Range("A13").Value = 107357705
$ws.Range("M13").Value = "äldre spår"
$ws.Range("Q13").Value = 421382.2179690556
$ws.Range("R13").Value = 7050382.422117869

# --- Row 14 ---
$ws.Range("A14").Value = 107357714
$ws.Range("M14").Value = "färska spår"
$ws.Range("Q14").Value = 421714.3643756237
$ws.Range("R14").Value = 7050134.446237801

# --- Row 15 ---
$ws.Range("A15").Value = 107357715
$ws.Range("M15").Value = "äldre spår"
$ws.Range("Q15").Value = 421720.3381778797
$ws.Range("R15").Value = 7050141.430032426

# --- Row 16 ---
$ws.Range("A16").Value = 107357710
$ws.Range("Q16").Value = 422003.2174995458
$ws.Range("R16").Value = 7049957.083066907

# --- Row 17 ---
$ws.Range("A17").Value = 107357713
$ws.Range("M17").Value = "färska spår"
$ws.Range("Q17").Value = 421707.8525994752
$ws.Range("R17").Value = 7050105.63680503

# --- Row 18 ---
$ws.Range("A18").Value = 107357716
$ws.Range("Q18").Value = 421633.682578478
$ws.Range("R18").Value = 7050192.592104986

# --- Row 19 (Tretåig hackspett -> Spillkråka) ---
$ws.Range("A19").Value = 105221881
$ws.Range("B19").Value = 56411
$ws.Range("E19").Value = 100049
$ws.Range("F19").Value = "Spillkråka"
$ws.Range("G19").Value = "Dryocopus martius"
$ws.Range("M19").Value = "äldre spår"
$ws.Range("P19").Value = "Djupsjö NV, Jmt"
$ws.Range("Q19").Value = 422185.3756590126
$ws.Range("R19").Value = 7049959.747113227
$ws.Range("S19").Value = 25
$ws.Range("AC19").Value = $null
$ws.Range("AW19").Value = "Johan Råghall"
$ws.Range("AX19").Value = "Johan Råghall, Benny Öwre, Lisa Gideonsson"

# --- Row 20 (Tretåig hackspett -> Spillkråka) ---
$ws.Range("A20").Value = 105221887
$ws.Range("B20").Value = 56411
$ws.Range("E20").Value = 100049
$ws.Range("F20").Value = "Spillkråka"
$ws.Range("G20").Value = "Dryocopus martius"
$ws.Range("P20").Value = "Djupsjö NV, Jmt"
$ws.Range("Q20").Value = 421493.4590253821
$ws.Range("R20").Value = 7050441.61621454
$ws.Range("S20").Value = 25
$ws.Range("AC20").Value = $null
$ws.Range("AW20").Value = "Johan Råghall"
$ws.Range("AX20").Value = "Johan Råghall, Benny Öwre, Lisa Gideonsson"

# --- Row 21 (Spillkråka -> Tretåig hackspett) ---
$ws.Range("A21").Value = 107357708
$ws.Range("B21").Value = 56395
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "1"
$ws.Range("M21").Value = "födosökande"
$ws.Range("Q21").Value = 422387.2765512859
$ws.Range("R21").Value = 7049729.271068533
$ws.Range("AC21").Value = $null

# --- Row 22 (Spillkråka -> Tretåig hackspett) ---
$ws.Range("A22").Value = 107357709
$ws.Range("B22").Value = 56395
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("P22").Value = "Djupsjö, Jmt"
$ws.Range("Q22").Value = 422308.8881526825
$ws.Range("R22").Value = 7049736.09021227
$ws.Range("S22").Value = 10
$ws.Range("AC22").Value = "ringhack"
$ws.Range("AW22").Value = "Benny Öwre"
$ws.Range("AX22").Value = "Benny Öwre, Johan Råghall, Lisa Gideonsson"

# --- Row 23 ---
$ws.Range("A23").Value = 107357720
$ws.Range("M23").Value = $null
$ws.Range("P23").Value = "Djupsjö, Jmt"
$ws.Range("Q23").Value = 422304.9171123825
$ws.Range("R23").Value = 7049737.970285944
$ws.Range("S23").Value = 10
$ws.Range("AC23").Value = "hack"
$ws.Range("AW23").Value = "Benny Öwre"
$ws.Range("AX23").Value = "Benny Öwre, Johan Råghall, Lisa Gideonsson"
